# Merge the "Search" + "Box" runs in the "Rectangle 11" shape (slide 1)
# into a single run reading "SearchBox", as shown in the target diff.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(41)
$tr = $sh.TextFrame.TextRange

# Remove the text of the second run ("Box") so PowerPoint collapses the
# paragraph down to the first run only (keeping its original rPr).
$tail = $tr.Characters(7, 3)
$tail.Text = ""

# Re-append "Box" onto the remaining run; since it's now a pure append to
# the existing single run, it is folded into that run's text rather than
# creating a new run, yielding one <a:r> containing "SearchBox".
$tr.Text = $tr.Text + "Box"
